$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (insert ownTeam/oppTeam columns before batsman, shifting columns right)
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# Existing row 3 data (Dubai DSC match) keeps its venue/date/result but gains ownTeam/oppTeam
# and shifts batsman stat columns right
$ws.Range("D3").Value = "Sunrisers Hyderabad"
$ws.Range("E3").Value = "Delhi Capitals"
$ws.Range("F3").Value = "Wriddhiman Saha " + [char]0x2020
$ws.Range("G3").Value = "87"
$ws.Range("H3").Value = "45"
$ws.Range("I3").Value = "12"
$ws.Range("J3").Value = "2"
$ws.Range("K3").Value = "193.33"

# New row 2: Sharjah vs Mumbai Indians match
$ws.Range("A2").Value = " Sharjah"
$ws.Range("B2").Value = " November 03 2020"
$ws.Range("C2").Value = "Sunrisers won by 10 wickets (with 17 balls remaining)"
$ws.Range("D2").Value = "Sunrisers Hyderabad"
$ws.Range("E2").Value = "Mumbai Indians"
$ws.Range("F2").Value = "Wriddhiman Saha " + [char]0x2020
$ws.Range("G2").Value = "58"
$ws.Range("H2").Value = "45"
$ws.Range("I2").Value = "7"
$ws.Range("J2").Value = "1"
$ws.Range("K2").Value = "128.88"

# New row 4: Sharjah vs Royal Challengers Bangalore match
$ws.Range("A4").Value = " Sharjah"
$ws.Range("B4").Value = " October 31 2020"
$ws.Range("C4").Value = "Sunrisers won by 5 wickets (with 35 balls remaining)"
$ws.Range("D4").Value = "Sunrisers Hyderabad"
$ws.Range("E4").Value = "Royal Challengers Bangalore"
$ws.Range("F4").Value = "Wriddhiman Saha " + [char]0x2020
$ws.Range("G4").Value = "39"
$ws.Range("H4").Value = "32"
$ws.Range("I4").Value = "4"
$ws.Range("J4").Value = "1"
$ws.Range("K4").Value = "121.87"

# New row 5: Abu Dhabi vs Kolkata Knight Riders match (the original row 2 content, moved down)
$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " September 26 2020"
$ws.Range("C5").Value = "KKR won by 7 wickets (with 12 balls remaining)"
$ws.Range("D5").Value = "Sunrisers Hyderabad"
$ws.Range("E5").Value = "Kolkata Knight Riders"
$ws.Range("F5").Value = "Wriddhiman Saha " + [char]0x2020
$ws.Range("G5").Value = "30"
$ws.Range("H5").Value = "31"
$ws.Range("I5").Value = "1"
$ws.Range("J5").Value = "1"
$ws.Range("K5").Value = "96.77"
